# Shift the recorded start/end timestamp strings in the Gantt "Main" sheet
# forward by a fixed offset (8 days, 3:43:28), matching the regenerated
# test fixture used by the commit "add spring data jpa and h2 repositories".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Map of old timestamp text -> new timestamp text, taken from the diff.
$map = @{
    "2019/04/13 19:51:19" = "2019/04/21 23:34:47"
    "2019/04/13 20:11:19" = "2019/04/21 23:54:47"
    "2019/04/13 20:31:19" = "2019/04/22 00:14:47"
    "2019/04/13 20:51:19" = "2019/04/22 00:34:47"
    "2019/04/13 21:11:19" = "2019/04/22 00:54:47"
    "2019/04/13 21:31:19" = "2019/04/22 01:14:47"
    "2019/04/13 21:51:19" = "2019/04/22 01:34:47"
    "2019/04/13 22:11:19" = "2019/04/22 01:54:47"
    "2019/04/13 22:31:19" = "2019/04/22 02:14:47"
    "2019/04/13 20:06:20" = "2019/04/21 23:49:48"
    "2019/04/13 20:21:21" = "2019/04/22 00:04:49"
    "2019/04/13 20:36:22" = "2019/04/22 00:19:50"
    "2019/04/13 20:51:23" = "2019/04/22 00:34:52"
    "2019/04/13 20:06:21" = "2019/04/21 23:49:49"
    "2019/04/13 20:21:24" = "2019/04/22 00:04:52"
    "2019/04/14 03:06:19" = "2019/04/22 06:49:47"
    "2019/04/14 10:21:19" = "2019/04/22 14:04:47"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# The timestamp values only ever live in columns D (start time) and F (end time).
# Walk row by row, column D before column F, matching the natural top-to-bottom,
# left-to-right reading order of the sheet.
for ($row = 1; $row -le $lastRow; $row++) {
    foreach ($col in @(4, 6)) {
        $cell = $ws.Cells.Item($row, $col)
        $val = $cell.Value2
        if ($val -ne $null -and $map.ContainsKey($val)) {
            $cell.Value2 = $map[$val]
        }
    }
}
